$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets("展览")
$sheet1.Range("F4").Value = 5981
$sheet1.Range("C5").Value = "北京·万游引力国潮动漫嘉年华s6"
$sheet1.Range("E5").Value = "2024.03.23 10:00-03.24 17:00"
$sheet1.Range("F5").Value = 3039
$sheet1.Range("G5").Value = 80
$sheet1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=79322"
$sheet1.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202402/wqACkjUk1708236212668.jpeg"
$sheet1.Range("C6").Value = "北京·排球少年ONLY"
$sheet1.Range("D6").Value = "永外高庄138号  大红门会展中心"
$sheet1.Range("E6").Value = "2024.03.23 10:00-03.23 17:00"
$sheet1.Range("F6").Value = 1298
$sheet1.Range("G6").Value = 70
$sheet1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=80510"
$sheet1.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202401/wNTz3awE1704441972575.jpeg"
$sheet1.Range("B7").Value = "2024-03-24"
$sheet1.Range("C7").Value = "北京· 次元音浪Million Mix——音乐番ONLY"
$sheet1.Range("D7").Value = "学清路38号 金码大厦"
$sheet1.Range("E7").Value = "2024.03.24 10:30-03.24 18:00"
$sheet1.Range("F7").Value = 448
$sheet1.Range("G7").Value = 68
$sheet1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=81640"
$sheet1.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202402/YhILflVA1706779569395.png"
$sheet1.Range("C8").Value = "北京·万游引力S6 知名声优 黑芝烧  张喆内场见面签售会"
$sheet1.Range("D8").Value = "半截塔路53号首创郎园station西门 郎园station中央车站文化中心"
$sheet1.Range("E8").Value = "2024.03.24 10:00-03.24 17:00"
$sheet1.Range("F8").Value = 108
$sheet1.Range("G8").Value = 238
$sheet1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=81855"
$sheet1.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202402/PkquxYeU1708225105344.jpeg"
$sheet1.Range("B9").Value = "2024-03-29"
$sheet1.Range("C9").Value = "北京·2024图书市集春季场"
$sheet1.Range("D9").Value = "建国路郎家园6号 郎园Vintage"
$sheet1.Range("E9").Value = "2024.03.29 14:00-03.31 20:00"
$sheet1.Range("F9").Value = 59
$sheet1.Range("G9").Value = 35.1
$sheet1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=81984"
$sheet1.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202402/Zi09QvTC1708571966640.jpeg"
$sheet1.Range("B10").Value = "2024-03-30"
$sheet1.Range("C10").Value = "北京·万象汇免费展"
$sheet1.Range("D10").Value = "滨河路178号 北京密云万象汇"
$sheet1.Range("E10").Value = "2024.03.30 13:00-03.30 17:00"
$sheet1.Range("F10").Value = 42
$sheet1.Range("G10").Value = 35
$sheet1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=83160"
$sheet1.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202403/9PoGPamI1710906663490.png"
$sheet1.Range("C11").Value = "北京·梦游园3.0代号鸢周年庆Only"
$sheet1.Range("D11").Value = "北花园路1号 超级蜂巢"
$sheet1.Range("E11").Value = "2024.03.30 10:00-03.30 17:00"
$sheet1.Range("F11").Value = 769
$sheet1.Range("G11").Value = 98
$sheet1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=81584"
$sheet1.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202402/ASPwEB9W1706844758149.jpeg"
$sheet1.Range("F12").Value = 337
$sheet1.Range("F13").Value = 4484
$sheet1.Range("F14").Value = 4484
$sheet1.Range("F16").Value = 99
$sheet1.Range("F17").Value = 132
$sheet1.Range("F18").Value = 24
$sheet1.Range("F20").Value = 82
$sheet1.Range("F21").Value = 6972
$sheet1.Range("F23").Value = 117
$sheet1.Range("F24").Value = 481
$sheet1.Range("F25").Value = 1282
$sheet1.Range("F26").Value = 6268
$sheet1.Range("F27").Value = 1654
$sheet1.Range("F28").Value = 17
$sheet1.Range("F30").Value = 6051
$sheet1.Range("F31").Value = 120
$sheet1.Range("F33").Value = 103
$sheet1.Range("F35").Value = 439
$sheet1.Range("F36").Value = 6138
$sheet1.Range("F37").Value = 18
$sheet1.Range("F39").Value = 89
$sheet1.Range("F41").Value = 9
$sheet1.Range("F42").Value = 2424
$sheet1.Range("F47").Value = 369
$sheet1.Range("F50").Value = 1044

$sheet2 = $wb.Worksheets("演出")
$sheet2.Range("F3").Value = 210
$sheet2.Range("F7").Value = 32
$sheet2.Range("F8").Value = 4

$sheet4 = $wb.Worksheets("全部类型")
$sheet4.Range("F3").Value = 5981
$sheet4.Range("F4").Value = 5981
$sheet4.Range("F5").Value = 3039
$sheet4.Range("F6").Value = 1298
$sheet4.Range("F7").Value = 448
$sheet4.Range("F9").Value = 210
$sheet4.Range("F10").Value = 42
$sheet4.Range("F12").Value = 337
$sheet4.Range("F13").Value = 4484
$sheet4.Range("F14").Value = 4484
$sheet4.Range("F16").Value = 99
$sheet4.Range("F17").Value = 132
$sheet4.Range("F18").Value = 24
$sheet4.Range("F20").Value = 82
$sheet4.Range("F21").Value = 6972
$sheet4.Range("F23").Value = 117
$sheet4.Range("F24").Value = 481
$sheet4.Range("F25").Value = 1282
$sheet4.Range("F27").Value = 6268
$sheet4.Range("F28").Value = 1654
$sheet4.Range("F30").Value = 4
$sheet4.Range("F32").Value = 6051
$sheet4.Range("F33").Value = 120
$sheet4.Range("F35").Value = 103
$sheet4.Range("F37").Value = 439
$sheet4.Range("F38").Value = 6138
$sheet4.Range("F39").Value = 18
$sheet4.Range("F41").Value = 89
$sheet4.Range("F42").Value = 9
$sheet4.Range("F44").Value = 2424
$sheet4.Range("F48").Value = 370
$sheet4.Range("F51").Value = 1044
